# Refresh the cryptocurrency price/volume snapshot (GitHub Actions bot update).
# Column D (Price) and E (Volume 1h) are plain text cells (not real numbers),
# so numeric-looking prices are forced back to text (NumberFormat "@" + Style
# "Normal" afterwards) to avoid Excel silently re-typing them as numbers while
# keeping the cells' style index identical to the untouched original cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.931.97"
$ws.Range("E2").Value = "  +0.40%  "

$ws.Range("D3").Value = "1.892.76"
$ws.Range("E3").Value = "  -0.33%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.8280"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +8.19%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "241.42"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.42%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.001"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3220"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +5.85%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "26.45"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.65%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07006"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.44%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08032"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.63%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7455"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.53%  "

$ws.Range("D13").Value = "1.894.67"
$ws.Range("E13").Value = "  -0.26%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.192"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.57%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "92.16"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.26%  "

$ws.Range("D16").Value = "29.943.48"
$ws.Range("E16").Value = "  +0.37%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.01"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.20%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.890"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.19%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "244.04"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.31%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007749"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.72%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.001"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.07%  "

$ws.Range("D22").Value = "2.149.10"
$ws.Range("E22").Value = "  +0.15%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.002"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.01%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.895"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.38%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1596"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +24.73%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "167.70"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.46%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.168"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.66%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.82"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.92%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.072"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.35%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.368"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.34%  "

$ws.Range("E31").Value = "  +0.50%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.254"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.20%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05641"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +7.42%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.072"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.26%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.270"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.44%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7307"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.83%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.721"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.07%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01909"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.04%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.780"
$ws.Range("D39").Style = "Normal"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.4403"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.09%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "71.87"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.01%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.943"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.15%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8419"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.78%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.001"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.01%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.886"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.41%  "

$ws.Range("B46").Value = "Aptos"
$ws.Range("C46").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.571"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.04%  "

$ws.Range("B47").Value = "Quant"
$ws.Range("C47").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "100.67"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.82%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.706"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.39%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "987.35"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +8.88%  "

$ws.Range("D50").Value = "2.045.07"
$ws.Range("E50").Value = "  +0.07%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "35.96"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.45%  "
